$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
    "B1" = "HSD-SMS - Movement for Self-Governing Democracy-Society for Moravia and Silesia (Hnutí za samosprávnou demokracii-Společnost pro Moravu a Slezsko, HSD-SMS)"
    "C1" = "KDU-ČSL - Christian and Democratic Union-Czech People's Party (Křesťanská a demokratická unie-Československá strana lidová, KDU-ČSL)"
    "D1" = "LB - Left Bloc Coalition (Koalice Levý blok, LB)"
    "E1" = "LSU - Liberal Social Union (Liberálně sociální unie, LSU)"
    "F1" = "ODA - Civic Democratic Alliance (Občanská demokratická aliance, ODA)"
    "G1" = "ODS - Civic Democratic Party (Občanská demokratická strana, ODS)"
    "H1" = "SPR-RSČ - Association for the Republic-Republican Party of Czechoslovakia (Sdružení pro republiku-Republikánská strana Československa, SPR-RSČ)"
    "I1" = "ČSSD - Czech Social Democratic Party (Česká strana sociálně demokratická, ČSSD), known until  as Czechoslovak Social Democracy (CSSD, Československá sociální demokracie)"
    "J1" = "KSČM - Communist Party of Bohemia and Moravia (Komunistická strana Čech a Moravy, KSČM)"
    "K1" = "US - Freedom Union (Unie svobodyUnie svobody-Demokratická unie, US)"
    "L1" = "SZ - Party of Greens (Strana zelených, SZ)"
    "M1" = "TOP09 - TOP09 (Tradition, Responsibility, Prosperity 09) (TOP09 (Tradice Odpovědnost Prosperita 09), TOP09)"
    "N1" = "VV - Public Affairs (Věci veřejné, VV)"
    "O1" = "ANO 2011 - Akce nespokojených občanů 2011 (Action of Dissatisfied Citizens 2011, ANO 2011)"
    "P1" = "Usvit - Dawn of Direct Democracy of Tomio Okamura (Úsvit Přímé Demokracie Tomia Okamury, Usvit)"
    "Q1" = "Piráti - Czech Pirate Party (Česká pirátská strana, Piráti)"
    "R1" = "SPD - Freedom and Direct Democracy (Svoboda a přímá demokracie, SPD)"
    "S1" = "STAN - Mayors and Independents (Starostové a Nezávislí, STAN)"
    "T1" = "PaS - Pirates and Mayors (Piráti a Starostové, PaS)"
    "U1" = "S - Together-ODS, KDU-ČSL, TOP 09 (Spolu-ODS, KDU-ČSL, TOP 09, S)"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

$ws.Range("I8").Value = 5
